$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.382.16"
$ws.Range("E2").Value = "  +0.05%  "
Set-TextValue $ws.Range("D3") "1.567.94"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.08%  "
Set-TextValue $ws.Range("D6") "291.57"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +2.46%  "
Set-TextValue $ws.Range("D8") "49.46"
$ws.Range("E8").Value = "  +0.36%  "
Set-TextValue $ws.Range("D9") "0.3408"
$ws.Range("E9").Value = "  +0.76%  "
Set-TextValue $ws.Range("D10") "0.07615"
$ws.Range("E10").Value = "  -0.07%  "
Set-TextValue $ws.Range("D11") "1.142"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  +0.04%  "
Set-TextValue $ws.Range("D13") "21.05"
$ws.Range("E13").Value = "  -0.81%  "
Set-TextValue $ws.Range("D14") "5.985"
$ws.Range("E14").Value = "  -1.17%  "
Set-TextValue $ws.Range("D15") "6.952"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D16") "0.00001134"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "1.547.39"
$ws.Range("E17").Value = "  -1.13%  "
Set-TextValue $ws.Range("D18") "90.03"
$ws.Range("E18").Value = "  +0.56%  "
Set-TextValue $ws.Range("D19") "0.06734"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.07%  "
Set-TextValue $ws.Range("D21") "16.60"
$ws.Range("E21").Value = "  +0.61%  "
Set-TextValue $ws.Range("D22") "6.203"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.28%  "
Set-TextValue $ws.Range("D24") "22.371.27"
$ws.Range("E24").Value = "  -0.03%  "
Set-TextValue $ws.Range("D25") "2.405"
$ws.Range("E25").Value = "  +1.29%  "
Set-TextValue $ws.Range("D26") "2.713"
$ws.Range("E26").Value = "  -6.60%  "
Set-TextValue $ws.Range("D27") "20.11"
$ws.Range("E27").Value = "  +0.51%  "
Set-TextValue $ws.Range("D28") "147.49"
$ws.Range("E28").Value = "  +1.49%  "
Set-TextValue $ws.Range("D29") "5.033"
$ws.Range("E29").Value = "  +1.15%  "
Set-TextValue $ws.Range("D30") "126.41"
$ws.Range("E30").Value = "  +0.67%  "
Set-TextValue $ws.Range("D31") "1.742.04"
$ws.Range("E31").Value = "  -0.20%  "
Set-TextValue $ws.Range("D32") "2.015"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "6.100"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D34") "0.9907"
$ws.Range("E34").Value = "  -5.32%  "
Set-TextValue $ws.Range("D35") "10.18"
$ws.Range("E35").Value = "  -0.24%  "
Set-TextValue $ws.Range("D36") "1.420"
$ws.Range("E36").Value = "  +9.12%  "
Set-TextValue $ws.Range("D37") "0.08478"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  -0.61%  "
Set-TextValue $ws.Range("D39") "0.2300"
$ws.Range("E39").Value = "  -0.66%  "
Set-TextValue $ws.Range("D40") "0.06478"
Set-TextValue $ws.Range("D41") "5.411"
Set-TextValue $ws.Range("D42") "11.38"
$ws.Range("E42").Value = "  -2.58%  "
Set-TextValue $ws.Range("D43") "0.6331"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D44") "1.001"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "14.04"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D46") "3.812"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5933"
$ws.Range("E47").Value = "  -0.66%  "
Set-TextValue $ws.Range("D48") "2.081"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D49") "1.275"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D50") "124.49"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.07314"
$ws.Range("E51").Value = "  +0.54%  "
